$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-D are plain text in the source data (coin names, links, prices
# formatted like "1.960.98"). Force text format first so Excel does not
# auto-coerce price-like strings (e.g. "230.80", "46.00") into numbers,
# which would silently drop the trailing/formatting zeros.
$ws.Range("B2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "36.597.39"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3
$ws.Range("D3").Value = "1.960.98"
$ws.Range("E3").Value = "  +0.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "244.83"
$ws.Range("E5").Value = "  +1.37%  "

# Row 6
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  +0.63%  "

# Row 7
$ws.Range("D7").Value = "59.11"
$ws.Range("E7").Value = "  +3.12%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "0.369"
$ws.Range("E9").Value = "  +0.55%  "

# Row 10
$ws.Range("D10").Value = "0.0859"
$ws.Range("E10").Value = "  +9.80%  "

# Row 11
$ws.Range("E11").Value = "  +1.48%  "

# Row 12
$ws.Range("D12").Value = "22.81"
$ws.Range("E12").Value = "  +5.83%  "

# Row 13
$ws.Range("D13").Value = "0.834"
$ws.Range("E13").Value = "  -0.76%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.253.25"
$ws.Range("E14").Value = "  +0.31%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "13.81"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16
$ws.Range("D16").Value = "5.28"
$ws.Range("E16").Value = "  -1.24%  "

# Row 17
$ws.Range("D17").Value = "1.963.88"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18
$ws.Range("D18").Value = "36.583.24"
$ws.Range("E18").Value = "  +1.83%  "

# Row 19
$ws.Range("D19").Value = "70.27"
$ws.Range("E19").Value = "  -0.54%  "

# Row 20
$ws.Range("E20").Value = "  +3.74%  "

# Row 21
$ws.Range("D21").Value = "230.80"
$ws.Range("E21").Value = "  -1.88%  "

# Row 22
$ws.Range("D22").Value = "5.10"
$ws.Range("E22").Value = "  -0.99%  "

# Row 23
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").Value = "2.49"
$ws.Range("E24").Value = "  -0.39%  "

# Row 25
$ws.Range("E25").Value = "  +2.20%  "

# Row 26
$ws.Range("E26").Value = "  -1.57%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "162.78"
$ws.Range("E27").Value = "  +1.42%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.137"
$ws.Range("E28").Value = "  +12.96%  "

# Row 29
$ws.Range("D29").Value = "19.64"
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("D30").Value = "0.118"
$ws.Range("E30").Value = "  +0.26%  "

# Row 31
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +6.80%  "

# Row 32
$ws.Range("E32").Value = "  -0.84%  "

# Row 33
$ws.Range("D33").Value = "0.0643"
$ws.Range("E33").Value = "  +5.57%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "4.34"
$ws.Range("E34").Value = "  -0.68%  "

# Row 35
$ws.Range("B35").Value = "THORChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D35").Value = "6.47"
$ws.Range("E35").Value = "  +6.80%  "

# Row 36
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("E37").Value = "  -1.79%  "

# Row 38
$ws.Range("E38").Value = "  -2.06%  "

# Row 39
$ws.Range("D39").Value = "3.08"
$ws.Range("E39").Value = "  +2.70%  "

# Row 40
$ws.Range("D40").Value = "0.0999"
$ws.Range("E40").Value = "  +1.82%  "

# Row 41
$ws.Range("E41").Value = "  -1.25%  "

# Row 42
$ws.Range("E42").Value = "  -0.05%  "

# Row 43
$ws.Range("D43").Value = "0.0211"
$ws.Range("E43").Value = "  +0.33%  "

# Row 44
$ws.Range("D44").Value = "16.44"
$ws.Range("E44").Value = "  +4.55%  "

# Row 45
$ws.Range("E45").Value = "  -2.70%  "

# Row 46
$ws.Range("D46").Value = "1.358.67"
$ws.Range("E46").Value = "  +2.32%  "

# Row 47
$ws.Range("D47").Value = "88.84"
$ws.Range("E47").Value = "  -2.15%  "

# Row 48
$ws.Range("E48").Value = "  -2.25%  "

# Row 49
$ws.Range("E49").Value = "  +0.47%  "

# Row 50
$ws.Range("D50").Value = "46.00"
$ws.Range("E50").Value = "  +5.54%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.142.31"
$ws.Range("E51").Value = "  +0.13%  "
